$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$kValues = @{
    2  = 3
    3  = 4
    4  = 8
    5  = 5
    6  = 11
    7  = 5
    8  = 7
    9  = 10
    10 = 7
    11 = 6
    12 = 5
    13 = 9
    14 = 5
    15 = 8
    16 = 5
    17 = 5
    18 = 8
    19 = 5
    20 = 6
    21 = 5
    22 = 6
    23 = 6
    24 = 8
    25 = 9
    26 = 6
    27 = 7
    28 = 2
    29 = 6
    30 = 5
    31 = 2
    32 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
